# "feat: experimenting with randomization"
#
# The stimuli/comprehension-questions sheet renames the distractor columns
# (distractor_1/2/3 -> distractor_a/b/c, and the accompanying
# distractor_1_span_text -> distractor_span_text) so the data no longer
# hard-codes a fixed distractor ordering, and clarifies the annotated-text
# header name (text_annotated_spans -> text_with_annotated_spans).
# Column positions / cell data are unchanged - only the header labels move.
#
# The shared-string table in the target workbook appends the new header
# strings in the order: distractor_a, text_with_annotated_spans,
# distractor_span_text, distractor_b, distractor_c - so we write the header
# cells in that same order (I1, E1, M1, J1, K1) to reproduce it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1").Value = "distractor_a"
$ws.Range("E1").Value = "text_with_annotated_spans"
$ws.Range("M1").Value = "distractor_span_text"
$ws.Range("J1").Value = "distractor_b"
$ws.Range("K1").Value = "distractor_c"

# Move the visible selection (scroll/selection state changed in the source
# commit from topLeftCell A4 / selection C6 to topLeftCell C1 / selection
# K1) - the engine round-trips the active selection via Range.Select().
$null = $ws.Range("K1").Select()
